$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.498.57"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "1.618.41"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.77"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").Value = "1.846.75"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("D13").Value = "1.621.48"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "27.483.87"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.71%  "
$ws.Range("E24").Value = "  +6.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("E26").Value = "  -1.26%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "1.445.38"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.28%  "
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.937"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.561"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.862"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.89%  "
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("E44").Value = "  -0.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.04%  "
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("D47").Value = "1.758.57"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("E51").Value = "  +1.54%  "
